$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellD = $ws.Range("D2")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "66.691.38"
$cellD.Style = $styleD
$ws.Range("E2").Value = "  +1.83%  "

$cellD = $ws.Range("D3")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "3.502.87"
$cellD.Style = $styleD
$ws.Range("E3").Value = "  +1.42%  "

$ws.Range("E4").Value = "  -0.12%  "

$cellD = $ws.Range("D5")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "593.67"
$cellD.Style = $styleD
$ws.Range("E5").Value = "  +2.42%  "

$cellD = $ws.Range("D6")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "169.38"
$cellD.Style = $styleD
$ws.Range("E6").Value = "  +0.92%  "

$ws.Range("E7").Value = "  -0.02%  "

$cellD = $ws.Range("D8")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "0.593"
$cellD.Style = $styleD
$ws.Range("E8").Value = "  +5.59%  "

$cellD = $ws.Range("D9")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "0.133"
$cellD.Style = $styleD
$ws.Range("E9").Value = "  +8.44%  "

$ws.Range("E10").Value = "  +1.48%  "

$ws.Range("E11").Value = "  +0.78%  "

$cellD = $ws.Range("D12")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "4.106.24"
$cellD.Style = $styleD
$ws.Range("E12").Value = "  +1.25%  "

$ws.Range("E14").Value = "  +2.88%  "

$cellD = $ws.Range("D15")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "0.0000182"
$cellD.Style = $styleD
$ws.Range("E15").Value = "  +4.01%  "

$cellD = $ws.Range("D16")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "66.708.82"
$cellD.Style = $styleD
$ws.Range("E16").Value = "  +1.94%  "

$cellD = $ws.Range("D17")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "3.491.54"
$cellD.Style = $styleD
$ws.Range("E17").Value = "  +0.74%  "

$ws.Range("E18").Value = "  +1.77%  "

$cellD = $ws.Range("D19")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "14.06"
$cellD.Style = $styleD
$ws.Range("E19").Value = "  +1.95%  "

$cellD = $ws.Range("D20")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "397.10"
$cellD.Style = $styleD
$ws.Range("E20").Value = "  +3.92%  "

$ws.Range("E21").Value = "  +0.67%  "

$cellD = $ws.Range("D22")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "73.40"
$cellD.Style = $styleD
$ws.Range("E22").Value = "  +2.41%  "

$ws.Range("E23").Value = "  +0.18%  "

$ws.Range("E24").Value = "  +2.91%  "

$ws.Range("E25").Value = "  +2.74%  "

$ws.Range("E26").Value = "  +2.54%  "

$cellD = $ws.Range("D27")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "0.181"
$cellD.Style = $styleD
$ws.Range("E27").Value = "  +0.40%  "

$cellD = $ws.Range("D28")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "0.999"
$cellD.Style = $styleD
$ws.Range("E28").Value = "  -0.06%  "

$cellD = $ws.Range("D29")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "6.32"
$cellD.Style = $styleD
$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("E30").Value = "  +1.28%  "

$ws.Range("E31").Value = "  +2.13%  "

$cellD = $ws.Range("D32")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "23.82"
$cellD.Style = $styleD
$ws.Range("E32").Value = "  +2.70%  "

$cellD = $ws.Range("D33")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "7.41"
$cellD.Style = $styleD
$ws.Range("E33").Value = "  +1.58%  "

$ws.Range("E34").Value = "  +6.20%  "

$cellD = $ws.Range("D35")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "162.30"
$cellD.Style = $styleD
$ws.Range("E35").Value = "  +1.28%  "

$ws.Range("E36").Value = "  +0.40%  "

$cellD = $ws.Range("D37")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "1.92"
$cellD.Style = $styleD
$ws.Range("E37").Value = "  +3.03%  "

$ws.Range("E38").Value = "  +2.82%  "

$cellD = $ws.Range("D39")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "4.69"
$cellD.Style = $styleD
$ws.Range("E39").Value = "  +5.34%  "

$cellD = $ws.Range("D40")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "0.0746"
$cellD.Style = $styleD
$ws.Range("E40").Value = "  +1.05%  "

$cellD = $ws.Range("D41")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "26.54"
$cellD.Style = $styleD
$ws.Range("E41").Value = "  +1.96%  "

$cellD = $ws.Range("D42")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "27.30"
$cellD.Style = $styleD
$ws.Range("E42").Value = "  +2.64%  "

$cellD = $ws.Range("D43")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "2.810.64"
$cellD.Style = $styleD
$ws.Range("E43").Value = "  -0.18%  "

$cellD = $ws.Range("D44")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "42.97"
$cellD.Style = $styleD
$ws.Range("E44").Value = "  -0.13%  "

$cellD = $ws.Range("D45")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "2.57"
$cellD.Style = $styleD
$ws.Range("E45").Value = "  +3.29%  "

$ws.Range("E46").Value = "  +2.10%  "

$cellD = $ws.Range("D47")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "342.79"
$cellD.Style = $styleD
$ws.Range("E47").Value = "  -1.10%  "

$cellD = $ws.Range("D48")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "1.10"
$cellD.Style = $styleD
$ws.Range("E48").Value = "  +2.45%  "

$cellD = $ws.Range("D49")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "34.16"
$cellD.Style = $styleD
$ws.Range("E49").Value = "  +5.58%  "

$cellD = $ws.Range("D50")
$styleD = $cellD.Style
$cellD.NumberFormat = "@"
$cellD.Value = "0.857"
$cellD.Style = $styleD
$ws.Range("E50").Value = "  +2.27%  "

$ws.Range("E51").Value = "  +2.30%  "
